$wb = $excel.ActiveWorkbook

# --- Sheets ---
# (GDPGR-alternate / GDPGR-bau recalc automatically via their formulas
#  that reference the Data sheet, so they don't need direct edits.)
$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")

# --- About sheet: update source text for September STEO / September 2020 ---
$wsAbout.Range("B6").Value = "January 2020 and September 2020"
$wsAbout.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of September 9,"

# --- Data sheet: update STEO label and revised GDP figures ---
$wsData.Range("A3").Value = "September STEO"
$wsData.Range("C3").Value = 18168
$wsData.Range("D3").Value = 18726

# --- Update selections to match final cursor positions ---
# Select Data's cell first, then About's, so "About" ends up as the
# active / front-most sheet (matching tabSelected="1" on About).
$wsData.Range("D4").Select() | Out-Null
$wsAbout.Range("A29").Select() | Out-Null

# --- Turn off iterative calculation (workbook no longer needs it) ---
$excel.Iteration = $false

# Recalculate everything so dependent formulas pick up new results
$excel.CalculateFullRebuild() | Out-Null
